$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcomes")

# Insert a new row above row 6 (shifts existing rows 6-12 down to 7-13)
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the UTI data
$ws.Range("A6").Value = 861
$ws.Range("B6").Value = "Earliest event of Urinary tract infections (UTI)"
$ws.Range("C6").Value = 9999

# Append two brand-new rows at the bottom (14 and 15)
$ws.Range("A14").Value = 1104
$ws.Range("B14").Value = "RBC Transfusion (adult relevant, no auto 1yr clean window)"
$ws.Range("C14").Value = 30

$ws.Range("A15").Value = 1105
$ws.Range("B15").Value = "Clostridium difficile - first episode"
$ws.Range("C15").Value = 9999
